# 19.xlsx — "Added support for longer quotes, fixed surplus numnber"
#
# Semantic changes applied:
#   1. Surplus/markup factor in column K reverted from 1.0565 to 1 on the
#      rows that still had the old value (K16, K17, K20, K23, K26, K29,
#      K30, K34) — "fixed surplus number".
#   2. Active cell / selection on the "Customer Quote" sheet moved from
#      A35 to G8 (a side effect of the editing session — "longer quotes"
#      work happened further up the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# --- 1. Surplus factor fix: 1.0565 -> 1 -------------------------------
$surplusRows = 16, 17, 20, 23, 26, 29, 30, 34
foreach ($row in $surplusRows) {
    $ws.Range("K$row").Value = 1
}

# --- 2. Restore the active selection left by the editing session ------
$ws.Range("G8").Select()
